# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de and
# zh-cn handback packages have been generated: the overall status moves
# from "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Target File" / "Latest Handback File" columns are populated for
# both languages, and the "Latest Handback DateTime" is stamped.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$sourceMd   = "49e75141-5af2-4885-aca5-2eac0a7b72b5.md"
$sourceUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdfb39a00a91c33dad652e1518ab65e37ce36439/e2e/49e75141-5af2-4885-aca5-2eac0a7b72b5.md"
$newStatus  = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2/F2 mirror the per-language Status cell on each sheet)
# ---------------------------------------------------------------------
$overview.Cells.Item(2, 5).Value = $newStatus   # Overview E2 (zh-cn column)
$overview.Cells.Item(2, 6).Value = $newStatus   # Overview F2 (de-de column)
$zhcn.Cells.Item(2, 3).Value = $newStatus        # zh-cn Status
$dede.Cells.Item(2, 3).Value = $newStatus        # de-de Status

# ---------------------------------------------------------------------
# 2. zh-cn row 2: Latest Target File / Latest Handback File / DateTime
# ---------------------------------------------------------------------
$zhcn.Hyperlinks.Add($zhcn.Cells.Item(2, 9), $sourceUrl, "", "", $sourceMd)
$zhcn.Cells.Item(2, 10).Value = "49e75141-5af2-4885-aca5-2eac0a7b72b5.c4f2c270380f0c69004b8fd0813c2898526f6860.zh-cn.xlf"
$zhcn.Cells.Item(2, 11).Value = "2016-09-05 05:05:16"

# ---------------------------------------------------------------------
# 3. de-de row 2: Latest Target File / Latest Handback File / DateTime
# ---------------------------------------------------------------------
$dede.Hyperlinks.Add($dede.Cells.Item(2, 9), $sourceUrl, "", "", $sourceMd)
$dede.Cells.Item(2, 10).Value = "49e75141-5af2-4885-aca5-2eac0a7b72b5.c4f2c270380f0c69004b8fd0813c2898526f6860.de-de.xlf"
$dede.Cells.Item(2, 11).Value = "2016-09-05 05:05:24"

# ---------------------------------------------------------------------
# 4. Column widths: widen the columns whose contents just grew so the
#    new text/links are readable (best-effort AutoFit-style resize).
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667

$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(9).ColumnWidth = 39.1666666666667
$zhcn.Columns.Item(10).ColumnWidth = 39.1666666666667

$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(9).ColumnWidth = 39.1666666666667
$dede.Columns.Item(10).ColumnWidth = 39.1666666666667
